$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.209789633750916
$ws.Range("B1").Value = 2.546462059020996
$ws.Range("C1").Value = 9.324790954589844
$ws.Range("D1").Value = 2.063595533370972
$ws.Range("E1").Value = 1.195932626724243
